$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1522.5
$ws.Range("I12").Value = 1676
$ws.Range("J12").Value = 1266.6666
$ws.Range("K12").Value = 1676
$ws.Range("L12").Value = 1266.6666
$ws.Range("M12").Value = -1506
$ws.Range("N12").Value = -1606.6666

$ws.Range("H28").Value = 301.5
$ws.Range("I28").Value = 360
$ws.Range("J28").Value = 172.8
$ws.Range("K28").Value = 360
$ws.Range("L28").Value = 172.8
$ws.Range("M28").Value = 125
$ws.Range("N28").Value = -1142.8

$ws.Range("H129").Value = 848.9778
$ws.Range("I129").Value = 499.4
$ws.Range("J129").Value = 892.675
$ws.Range("K129").Value = 1498.2
$ws.Range("L129").Value = 2678.025
$ws.Range("M129").Value = 3501.8
$ws.Range("N129").Value = -12678.025

$ws.Range("H137").Value = 37358.723
$ws.Range("I137").Value = 3664.2144
$ws.Range("J137").Value = 68806.92999999999
$ws.Range("K137").Value = 10992.6432
$ws.Range("L137").Value = 206420.79
$ws.Range("M137").Value = -8442.643199999999
$ws.Range("N137").Value = -211520.79

$ws.Range("H141").Value = 1622.5758
$ws.Range("I141").Value = 1204.3549
$ws.Range("J141").Value = 8105
$ws.Range("K141").Value = 3613.0647
$ws.Range("L141").Value = 24315
$ws.Range("M141").Value = 1566.9353
$ws.Range("N141").Value = -34675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1010.069
$ws.Range("I2").Value = 776.48834
$ws.Range("J2").Value = 1679.6666
$ws.Range("K2").Value = 776.48834
$ws.Range("L2").Value = 1679.6666
$ws.Range("M2").Value = -663.48834
$ws.Range("N2").Value = -1905.6666

$ws.Range("H3").Value = 350
$ws.Range("I3").Value = 350
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 350
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -235

$ws.Range("H32").Value = 21372.717
$ws.Range("I32").Value = 25405.465
$ws.Range("J32").Value = 4031.9
$ws.Range("K32").Value = 25405.465
$ws.Range("L32").Value = 4031.9
$ws.Range("M32").Value = -25118.465
$ws.Range("N32").Value = -4605.9

$ws.Range("H61").Value = 453407.34
$ws.Range("I61").Value = 722751.2
$ws.Range("J61").Value = 4500.933
$ws.Range("K61").Value = 722751.2
$ws.Range("L61").Value = 4500.933
$ws.Range("M61").Value = -722539.2
$ws.Range("N61").Value = -4924.933

$ws.Range("H63").Value = 2405838.2
$ws.Range("I63").Value = 2158
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2158
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -1472
$ws.Range("N63").Value = -31251372

$ws.Range("H66").Value = 2405838.2
$ws.Range("I66").Value = 2158
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 10790
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -7358
$ws.Range("N66").Value = -156256864

$ws.Range("H74").Value = 4363.933
$ws.Range("I74").Value = 6243.625
$ws.Range("J74").Value = 2215.7144
$ws.Range("K74").Value = 6243.625
$ws.Range("L74").Value = 2215.7144
$ws.Range("M74").Value = -5369.625
$ws.Range("N74").Value = -3963.7144

$ws.Range("H77").Value = 4363.933
$ws.Range("I77").Value = 6243.625
$ws.Range("J77").Value = 2215.7144
$ws.Range("K77").Value = 31218.125
$ws.Range("L77").Value = 11078.572
$ws.Range("M77").Value = -26850.125
$ws.Range("N77").Value = -19814.572

$ws.Range("H116").Value = 1010.069
$ws.Range("I116").Value = 776.48834
$ws.Range("J116").Value = 1679.6666
$ws.Range("K116").Value = 776.48834
$ws.Range("L116").Value = 1679.6666
$ws.Range("M116").Value = 1517.51166
$ws.Range("N116").Value = -6267.6666

$ws.Range("H132").Value = 22932.916
$ws.Range("I132").Value = 1832.7778
$ws.Range("J132").Value = 86233.336
$ws.Range("K132").Value = 5498.3334
$ws.Range("L132").Value = 258700.008
$ws.Range("M132").Value = -2968.3334
$ws.Range("N132").Value = -263760.008

$ws.Range("H136").Value = 453407.34
$ws.Range("I136").Value = 722751.2
$ws.Range("J136").Value = 4500.933
$ws.Range("K136").Value = 2168253.6
$ws.Range("L136").Value = 13502.799
$ws.Range("M136").Value = -2165703.6
$ws.Range("N136").Value = -18602.799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1010.069
$ws.Range("I3").Value = 776.48834
$ws.Range("J3").Value = 1679.6666
$ws.Range("K3").Value = 776.48834
$ws.Range("L3").Value = 1679.6666
$ws.Range("M3").Value = -662.48834
$ws.Range("N3").Value = -1907.6666

$ws.Range("H35").Value = 36500
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 36500
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 36500
$ws.Range("N35").Value = -37120

$ws.Range("H82").Value = 13164.75
$ws.Range("I82").Value = 6244.2
$ws.Range("J82").Value = 47767.5
$ws.Range("K82").Value = 6244.2
$ws.Range("L82").Value = 47767.5
$ws.Range("M82").Value = -5861.2
$ws.Range("N82").Value = -48533.5

$ws.Range("H85").Value = 13164.75
$ws.Range("I85").Value = 6244.2
$ws.Range("J85").Value = 47767.5
$ws.Range("K85").Value = 6244.2
$ws.Range("L85").Value = 47767.5
$ws.Range("M85").Value = -4918.2
$ws.Range("N85").Value = -50419.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 109
$ws.Range("N2").ClearContents()

$ws.Range("H31").Value = 9319.268
$ws.Range("I31").Value = 15013.5
$ws.Range("J31").Value = 2749
$ws.Range("K31").Value = 15013.5
$ws.Range("L31").Value = 2749
$ws.Range("M31").Value = -14718.5
$ws.Range("N31").Value = -3339

$ws.Range("H34").Value = 9319.268
$ws.Range("I34").Value = 15013.5
$ws.Range("J34").Value = 2749
$ws.Range("K34").Value = 15013.5
$ws.Range("L34").Value = 2749
$ws.Range("M34").Value = -14811.5
$ws.Range("N34").Value = -3153

$ws.Range("H58").Value = 30975.883
$ws.Range("I58").Value = 1545
$ws.Range("J58").Value = 101610
$ws.Range("K58").Value = 1545
$ws.Range("L58").Value = 101610
$ws.Range("M58").Value = -1342
$ws.Range("N58").Value = -102016

$ws.Range("H136").Value = 30975.883
$ws.Range("I136").Value = 1545
$ws.Range("J136").Value = 101610
$ws.Range("K136").Value = 4635
$ws.Range("L136").Value = 304830
$ws.Range("M136").Value = -2085
$ws.Range("N136").Value = -309930

$ws.Range("H141").Value = 65000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 65000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 821.7941
$ws.Range("I5").Value = 794.6667
$ws.Range("J5").Value = 831.5599999999999
$ws.Range("K5").Value = 2384.0001
$ws.Range("L5").Value = 2494.68
$ws.Range("M5").Value = -2272.0001
$ws.Range("N5").Value = -2718.68

$ws.Range("H7").Value = 15
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 45
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 67

$ws.Range("H56").Value = 6323
$ws.Range("I56").Value = 6323
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 6323
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -5793

$ws.Range("H68").Value = 3724.1316
$ws.Range("I68").Value = 1100
$ws.Range("J68").Value = 3869.9167
$ws.Range("K68").Value = 3300
$ws.Range("L68").Value = 11609.7501
$ws.Range("M68").Value = -2489
$ws.Range("N68").Value = -13231.7501

$ws.Range("H71").Value = 3724.1316
$ws.Range("I71").Value = 1100
$ws.Range("J71").Value = 3869.9167
$ws.Range("K71").Value = 9900
$ws.Range("L71").Value = 34829.2503
$ws.Range("M71").Value = -5844
$ws.Range("N71").Value = -42941.2503

$ws.Range("H80").Value = 2966.6667
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2966.6667
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 8900.000100000001
$ws.Range("N80").Value = -10772.0001

$ws.Range("H81").Value = 4015.6667
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 4015.6667
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 12047.0001
$ws.Range("N81").Value = -14293.0001

$ws.Range("H83").Value = 2966.6667
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2966.6667
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 26700.0003
$ws.Range("N83").Value = -36060.0003

$ws.Range("H84").Value = 4015.6667
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 4015.6667
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 36141.0003
$ws.Range("N84").Value = -47373.0003

$ws.Range("H107").Value = 4275.433
$ws.Range("I107").Value = 25449
$ws.Range("J107").Value = 1017.96155
$ws.Range("K107").Value = 76347
$ws.Range("L107").Value = 3053.88465
$ws.Range("M107").Value = -74427
$ws.Range("N107").Value = -6893.88465

$ws.Range("H121").Value = 4022.8708
$ws.Range("I121").Value = 516
$ws.Range("J121").Value = 4697.269
$ws.Range("K121").Value = 1548
$ws.Range("L121").Value = 14091.807
$ws.Range("M121").Value = -238
$ws.Range("N121").Value = -16711.807

$ws.Range("H131").Value = 105009.09
$ws.Range("I131").Value = 794.2857
$ws.Range("J131").Value = 113205.766
$ws.Range("K131").Value = 2382.8571
$ws.Range("L131").Value = 339617.298
$ws.Range("M131").Value = 2657.1429
$ws.Range("N131").Value = -349697.298

$ws.Range("H135").Value = 821.7941
$ws.Range("I135").Value = 794.6667
$ws.Range("J135").Value = 831.5599999999999
$ws.Range("K135").Value = 7152.0003
$ws.Range("L135").Value = 7484.039999999999
$ws.Range("M135").Value = -4617.0003
$ws.Range("N135").Value = -12554.04

$ws.Range("H139").Value = 1795.2307
$ws.Range("I139").Value = 1456.9166
$ws.Range("J139").Value = 5855
$ws.Range("K139").Value = 4370.7498
$ws.Range("L139").Value = 17565
$ws.Range("M139").Value = 769.2502000000004
$ws.Range("N139").Value = -27845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2416.6
$ws.Range("I22").Value = 2800.25
$ws.Range("J22").Value = 882
$ws.Range("K22").Value = 2800.25
$ws.Range("L22").Value = 882
$ws.Range("M22").Value = -2505.25
$ws.Range("N22").Value = -1472

$ws.Range("H27").Value = 2416.6
$ws.Range("I27").Value = 2800.25
$ws.Range("J27").Value = 882
$ws.Range("K27").Value = 2800.25
$ws.Range("L27").Value = 882
$ws.Range("M27").Value = -2693.25
$ws.Range("N27").Value = -1096

$ws.Range("H68").Value = 3581.7222
$ws.Range("I68").Value = 1718.7778
$ws.Range("J68").Value = 5444.6665
$ws.Range("K68").Value = 1718.7778
$ws.Range("L68").Value = 5444.6665
$ws.Range("M68").Value = -969.7778000000001
$ws.Range("N68").Value = -6942.6665

$ws.Range("H71").Value = 3581.7222
$ws.Range("I71").Value = 1718.7778
$ws.Range("J71").Value = 5444.6665
$ws.Range("K71").Value = 8593.889000000001
$ws.Range("L71").Value = 27223.3325
$ws.Range("M71").Value = -4849.889000000001
$ws.Range("N71").Value = -34711.3325

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 25000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 25000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

$ws.Range("H132").Value = 3445.2727
$ws.Range("I132").Value = 3050
$ws.Range("J132").Value = 4499.3335
$ws.Range("K132").Value = 9150
$ws.Range("L132").Value = 13498.0005
$ws.Range("M132").Value = -6620
$ws.Range("N132").Value = -18558.0005

$ws.Range("H139").Value = 52715
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 52715
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 52715
$ws.Range("N139").Value = -62995
